$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row for the new "2509" period, right before the empty
#        gap / footer rows (old row 19), shifting the footer rows (old 23,
#        24) down to (24, 25). ---
$ws.Rows("19").Insert()

# --- 2. Re-apply table borders/format so the table looks contiguous again:
#        the row that used to be the last data row (18) had the special
#        "bottom border" styling; now that a new last row (19) has been
#        added, row 18 should look like a normal middle row (like 16/17)
#        and the new row 19 should inherit the "bottom border" styling
#        that row 18 used to have. ---
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)

$ws.Range("B17:J17").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3. Fill in the new worker/period row (B19:J19) with the same
#        worker data as the other rows, but period 2509. ---
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "32942478"
$ws.Range("D19").Value = "LILIAN VANESSA TORRES BORGE"
$ws.Range("E19").Value = "2509"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

# --- 4. Center-align the "Periodo Mora" column for all data rows. ---
$ws.Range("E16:E19").HorizontalAlignment = -4108

# --- 5. Update the summary figures: total overdue value and period count. ---
$ws.Range("E11").Value = 227760
$ws.Range("F13").Value = 4
